# Automatically update sets from .txt files
#
# - Clears the stale "N/A" placeholder out of D19 (Collection column for
#   set 40825) so it matches the other "unknown collection" rows.
# - Appends two newly scraped LEGO sets (10313, 10329) as rows 20-21,
#   mirroring the existing table layout (ID_Set, Nom_Set, nbPieces,
#   Collection, Image_URL, URL_Lego, URL_Auchan, URL_Leclerc,
#   URL_Carrefour, URL_AvenueDeLaBrique).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (matching the rest of the sheet,
# where even numeric-looking values like IDs/piece counts are stored as
# text) without leaving the cell permanently tagged with a text number
# format, so the cell's style stays "no explicit style" just like its
# neighbours.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- D19: "N/A" -> blank -----------------------------------------------
$ws.Cells.Item(19, 4).Value = ""

# --- Row 20: set 10313 ---------------------------------------------------
Set-TextValue $ws.Cells.Item(20, 1) "10313"
$ws.Cells.Item(20, 2).Value = "Bouquet de fleurs sauvages"
Set-TextValue $ws.Cells.Item(20, 3) "939"
$ws.Cells.Item(20, 4).Value = "The Botanical Collection"
$ws.Cells.Item(20, 5).Value = "https://www.lego.com/cdn/cs/set/assets/bltc4a6c2103a34f22e/10313_alt2.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1"
$ws.Cells.Item(20, 6).Value = "https://www.lego.com/fr-fr/product/10313"

# --- Row 21: set 10329 ---------------------------------------------------
Set-TextValue $ws.Cells.Item(21, 1) "10329"
$ws.Cells.Item(21, 2).Value = "Les plantes miniatures"
Set-TextValue $ws.Cells.Item(21, 3) "758"
$ws.Cells.Item(21, 4).Value = "The Botanical Collection"
$ws.Cells.Item(21, 5).Value = "https://www.lego.com/cdn/cs/set/assets/bltb2f845ffd52a25b0/10329.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1"
$ws.Cells.Item(21, 6).Value = "https://www.lego.com/fr-fr/product/10329"
